$d = $word.ActiveDocument

# The document body contains a single (empty) paragraph. The authored
# change sets that paragraph mark's language to English (US) - i.e. the
# paragraph's run properties gain <w:rPr><w:lang w:val="en-US"/></w:rPr>.
$d.Paragraphs(1).Range.LanguageID = "en-US"
